# Update the timestamped e-mail addresses on the "UsuariosRegistro" sheet.
# The diff replaces the old timestamp 20251112_211458 with 20251112_215226
# in each of the 5 e-mail values found in column C (E-Mail), rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$oldStamp = "20251112_211458"
$newStamp = "20251112_215226"

for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C = E-Mail
    $current = $cell.Value2
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}
